$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K (strikeout) values for rows 2-16 (column G)
$values = @{
    2  = 8
    3  = 10
    4  = 9
    5  = 9
    6  = 7
    7  = 9
    8  = 7
    9  = 7
    10 = 3
    11 = 6
    12 = 4
    13 = 9
    14 = 3
    15 = 5
    16 = 7
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
